{"js": "// Apply four text corrections (Korean marketing copy wording tweaks) to the\n// document body. Each old string is unique in the document, so a simple\n// search + full-match replace (location \"Replace\") is safe and preserves\n// the run's formatting (rPr) because insertText(\"Replace\") rewrites the\n// text of the matched range in place.\nconst replacements = [\n  {\n    find: \"\ubcf5\uc6d0\ub825\uc774 \uc6b0\uc218\ud55c \uace0\uae09 \ubcf4\uc548 \uc81c\ud488\uc778 Contoso CipherGuard Sentinel X7\uc740 \ucef4\ud4e8\ud130 \ub124\ud2b8\uc6cc\ud06c \uc778\ud504\ub77c\ub97c \uac15\ud654\ud558\uc5ec \ub2e4\uc591\ud55c \uc704\ud611\uacfc \ucde8\uc57d\uc131\uc744 \ubc29\uc9c0\ud560 \uc218 \uc788\ub3c4\ub85d \uc138\uc2ec\ud558\uac8c \uc81c\uc791\ub418\uc5c8\uc2b5\ub2c8\ub2e4. \",\n    replace: \"Contoso CipherGuard Sentinel X7\uc740 \ub2e4\uc591\ud55c \uc704\ud611\uacfc \ucde8\uc57d\uc131\uc5d0 \ub300\ube44\ud558\uc5ec \ucef4\ud4e8\ud130 \ub124\ud2b8\uc6cc\ud06c \uc778\ud504\ub77c\ub97c \uac15\ud654\ud558\ub3c4\ub85d \uc138\uc2ec\ud558\uac8c \uc124\uacc4\ub41c \ubcf5\uc6d0\ub825 \uc788\ub294 \uace0\uae09 \ubcf4\uc548 \uc81c\ud488\uc785\ub2c8\ub2e4. \",\n  },\n  {\n    find: \" \ub2e4\uacc4\uce35 \ubc29\uc5b4 \uc811\uadfc \ubc29\uc2dd\uc744 \uc0ac\uc6a9\ud558\ub294 \uc5d4\ub4dc\ud3ec\uc778\ud2b8 \ubcf4\uc548 \ubaa8\ub4c8\uc740 \ubc14\uc774\ub7ec\uc2a4 \ubc31\uc2e0, \ub9ec\uc6e8\uc5b4 \ubc29\uc9c0 \ubc0f \ud638\uc2a4\ud2b8 \uae30\ubc18 \uce68\uc785 \ubc29\uc9c0 \uae30\ub2a5\uc744 \ud1b5\ud569\ud569\ub2c8\ub2e4. \",\n    replace: \" \ub2e4\uc911 \uacc4\uce35 \ubc29\uc5b4 \uc811\uadfc \ubc29\uc2dd\uc744 \uc0ac\uc6a9\ud558\ub294 \uc5d4\ub4dc\ud3ec\uc778\ud2b8 \ubcf4\uc548 \ubaa8\ub4c8\uc740 \ubc14\uc774\ub7ec\uc2a4 \ubc31\uc2e0, \ub9ec\uc6e8\uc5b4 \ubc29\uc9c0 \ubc0f \ud638\uc2a4\ud2b8 \uae30\ubc18 \uce68\uc785 \ubc29\uc9c0 \uae30\ub2a5\uc744 \ud1b5\ud569\ud569\ub2c8\ub2e4. \",\n  },\n  {\n    find: \" Windows Server 2019 \uc774\uc0c1, CentOS 8 \ub610\ub294 \ud574\ub2f9 \ubc84\uc804\uacfc \ud638\ud658\",\n    replace: \" Windows Server 2019 \uc774\uc0c1, CentOS 8 \ub610\ub294 \ub3d9\uae09 \ubc84\uc804\uacfc \ud638\ud658\",\n  },\n  {\n    find: \" Contoso\ub294 Contoso CipherGuard Sentinel X7\uacfc \uad00\ub828\ub41c \uae30\uc220 \ubb38\uc81c \ub610\ub294 \ubb38\uc758\uc5d0 \ub300\ud55c \uc2e0\uc18d\ud55c \uc9c0\uc6d0\uc744 \ubcf4\uc7a5\ud558\uae30 \uc704\ud574 \uc804\uc6a9 24/7 \uc9c0\uc6d0 \ud300\uc744 \uc81c\uacf5\ud569\ub2c8\ub2e4.\",\n    replace: \" Contoso\ub294 Contoso CipherGuard Sentinel X7\uacfc \uad00\ub828\ub41c \uae30\uc220 \ubb38\uc81c \ub610\ub294 \ubb38\uc758\uc5d0 \ub300\ud55c \uc2e0\uc18d\ud55c \uc9c0\uc6d0\uc744 \ubcf4\uc7a5\ud558\uae30 \uc704\ud574 \uc5f0\uc911\ubb34\ud734 \uc0c1\uc2dc \uc804\ub2f4 \uc9c0\uc6d0\ud300\uc744 \uc81c\uacf5\ud569\ub2c8\ub2e4.\",\n  },\n];\n\nconst body = context.document.body;\n\nfor (const { find, replace } of replacements) {\n  const results = body.search(find, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(replace, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Apply four text corrections (Korean marketing copy wording tweaks) via\n# Word's Find/Replace. Each FindText string is unique in the document, so\n# a single Execute(..., Replace:=wdReplaceAll) per pair is sufficient and\n# keeps the host run's character formatting intact.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{\n        Find    = \"\ubcf5\uc6d0\ub825\uc774 \uc6b0\uc218\ud55c \uace0\uae09 \ubcf4\uc548 \uc81c\ud488\uc778 Contoso CipherGuard Sentinel X7\uc740 \ucef4\ud4e8\ud130 \ub124\ud2b8\uc6cc\ud06c \uc778\ud504\ub77c\ub97c \uac15\ud654\ud558\uc5ec \ub2e4\uc591\ud55c \uc704\ud611\uacfc \ucde8\uc57d\uc131\uc744 \ubc29\uc9c0\ud560 \uc218 \uc788\ub3c4\ub85d \uc138\uc2ec\ud558\uac8c \uc81c\uc791\ub418\uc5c8\uc2b5\ub2c8\ub2e4. \"\n        Replace = \"Contoso CipherGuard Sentinel X7\uc740 \ub2e4\uc591\ud55c \uc704\ud611\uacfc \ucde8\uc57d\uc131\uc5d0 \ub300\ube44\ud558\uc5ec \ucef4\ud4e8\ud130 \ub124\ud2b8\uc6cc\ud06c \uc778\ud504\ub77c\ub97c \uac15\ud654\ud558\ub3c4\ub85d \uc138\uc2ec\ud558\uac8c \uc124\uacc4\ub41c \ubcf5\uc6d0\ub825 \uc788\ub294 \uace0\uae09 \ubcf4\uc548 \uc81c\ud488\uc785\ub2c8\ub2e4. \"\n    },\n    @{\n        Find    = \" \ub2e4\uacc4\uce35 \ubc29\uc5b4 \uc811\uadfc \ubc29\uc2dd\uc744 \uc0ac\uc6a9\ud558\ub294 \uc5d4\ub4dc\ud3ec\uc778\ud2b8 \ubcf4\uc548 \ubaa8\ub4c8\uc740 \ubc14\uc774\ub7ec\uc2a4 \ubc31\uc2e0, \ub9ec\uc6e8\uc5b4 \ubc29\uc9c0 \ubc0f \ud638\uc2a4\ud2b8 \uae30\ubc18 \uce68\uc785 \ubc29\uc9c0 \uae30\ub2a5\uc744 \ud1b5\ud569\ud569\ub2c8\ub2e4. \"\n        Replace = \" \ub2e4\uc911 \uacc4\uce35 \ubc29\uc5b4 \uc811\uadfc \ubc29\uc2dd\uc744 \uc0ac\uc6a9\ud558\ub294 \uc5d4\ub4dc\ud3ec\uc778\ud2b8 \ubcf4\uc548 \ubaa8\ub4c8\uc740 \ubc14\uc774\ub7ec\uc2a4 \ubc31\uc2e0, \ub9ec\uc6e8\uc5b4 \ubc29\uc9c0 \ubc0f \ud638\uc2a4\ud2b8 \uae30\ubc18 \uce68\uc785 \ubc29\uc9c0 \uae30\ub2a5\uc744 \ud1b5\ud569\ud569\ub2c8\ub2e4. \"\n    },\n    @{\n        Find    = \" Windows Server 2019 \uc774\uc0c1, CentOS 8 \ub610\ub294 \ud574\ub2f9 \ubc84\uc804\uacfc \ud638\ud658\"\n        Replace = \" Windows Server 2019 \uc774\uc0c1, CentOS 8 \ub610\ub294 \ub3d9\uae09 \ubc84\uc804\uacfc \ud638\ud658\"\n    },\n    @{\n        Find    = \" Contoso\ub294 Contoso CipherGuard Sentinel X7\uacfc \uad00\ub828\ub41c \uae30\uc220 \ubb38\uc81c \ub610\ub294 \ubb38\uc758\uc5d0 \ub300\ud55c \uc2e0\uc18d\ud55c \uc9c0\uc6d0\uc744 \ubcf4\uc7a5\ud558\uae30 \uc704\ud574 \uc804\uc6a9 24/7 \uc9c0\uc6d0 \ud300\uc744 \uc81c\uacf5\ud569\ub2c8\ub2e4.\"\n        Replace = \" Contoso\ub294 Contoso CipherGuard Sentinel X7\uacfc \uad00\ub828\ub41c \uae30\uc220 \ubb38\uc81c \ub610\ub294 \ubb38\uc758\uc5d0 \ub300\ud55c \uc2e0\uc18d\ud55c \uc9c0\uc6d0\uc744 \ubcf4\uc7a5\ud558\uae30 \uc704\ud574 \uc5f0\uc911\ubb34\ud734 \uc0c1\uc2dc \uc804\ub2f4 \uc9c0\uc6d0\ud300\uc744 \uc81c\uacf5\ud569\ub2c8\ub2e4.\"\n    }\n)\n\nforeach ($pair in $replacements) {\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    # Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,\n    #         MatchAllWordForms, Forward, Wrap, Format, ReplaceWith, Replace)\n    $find.Execute($pair.Find, $true, $false, $false, $false, $false, $true, 1, $false, $pair.Replace, 2) | Out-Null\n}\n"}
